$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '308.12'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-4.27%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '53.96'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '9.43%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.083'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-4.97%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07844'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-2.62%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.531'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-1.68%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.381'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-1.29%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.771'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '7.97%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1239'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-2.91%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2015'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2.49%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.04714'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.13%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09391'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-2.22%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1044'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.10%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001263'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-4.99%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005776'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.26%'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '2,018.44%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.327'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.48%'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.64%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3418'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-2.52%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.005'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.08%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1364'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.69%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2915'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-5.64%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04161'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.84%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001262'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-4.00%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.003949'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-8.89%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001349'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.23%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02598'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-4.80%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05853'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-5.57%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.009853'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-8.99%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007940'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.51%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1439'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-1.54%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008194'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '3.84%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008358'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-3.24%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3362'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-3.84%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00007311'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '10.23%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000749'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.03%'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '3.20%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002617'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002098'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.03%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0001998'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.03%'
